# Update the title on slide 1 ("BigBlueButton" -> "VISAR-Qtune").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

$tr.Text = "VISAR-Qtune"
# The author's run carries an indeterminate/"none" run language (lang="")
# rather than the usual "en-US" - reproduce that via LanguageID.
$tr.LanguageID = ""
